$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 'Maximilian_De_Junious(2cn)'
$ws.Range("E3").Value = 'Ramon#Nunez_Gomez&3bn='
$ws.Range("E4").Value = 'Uelkue_Oemer)Uellaegoess)1an&'
$ws.Range("E5").Value = 'Isis_Lanpher!4cn,'
$ws.Range("E6").Value = 'Maximilian-Galvin_2an!'
$ws.Range("E7").Value = 'Jannette(Laspina!4cn#'
$ws.Range("E8").Value = 'Marg(Dodich)3cn)'
$ws.Range("E9").Value = 'Carisa.Bannowsky,2an.'
$ws.Range("E10").Value = 'David^Waisath_4cn,'
$ws.Range("E11").Value = 'Paulette=von_Reddig-Piette(2cn='
$ws.Range("E12").Value = 'Kirby(Latona!2an%'
$ws.Range("E13").Value = 'Reed#Homewood,2bn('
$ws.Range("E14").Value = 'Blair.Pallafor-Zedian-5cn_'
$ws.Range("E15").Value = 'Lon=Senemounnarat-Quillian#2cn%'
$ws.Range("E16").Value = 'Vada.Isaac-2bn_'
$ws.Range("E17").Value = 'Jeanett!Plancarte-4bn,'
$ws.Range("E18").Value = 'Alex,Berteotti-Stirn-l3hr3r%'
$ws.Range("E19").Value = 'Robyn!Strycker&3an%'
$ws.Range("E20").Value = 'Camille^Von_Verrill%3bn('
$ws.Range("E21").Value = 'Franz_Michael_Leopold)Deschner-4cn.'
$ws.Range("E22").Value = 'Veola=Franzi#l3hr3r='
$ws.Range("E23").Value = 'Chantelle-Cringle%3cn!'
$ws.Range("E24").Value = 'Britney#Kosh(5bn%'
$ws.Range("E25").Value = 'Clayton^Derouchie^4bn&'
$ws.Range("E26").Value = 'Beverlee%Doutt(5bn('
$ws.Range("E27").Value = 'Alma%Munley!1bn#'
$ws.Range("E28").Value = 'Thad,Dornbos-5an('
$ws.Range("E29").Value = 'Arvilla_Mahala&2cn&'
$ws.Range("E30").Value = 'Mirza)Ellingwood(l3hr3r!'
$ws.Range("E31").Value = 'Francie%de_Cardinalli-Sciola.3an#'
$ws.Range("E32").Value = 'IRENEE#Gundry^l3hr3r_'
$ws.Range("E33").Value = 'IRENEE,Pinedo(3cn-'
$ws.Range("E34").Value = 'Mirza-Pinsky!1bn!'
$ws.Range("E35").Value = 'Francie=Pinsky1&1bn-'
$ws.Range("E36").Value = 'Goldie=Pinsky2(1bn)'
$ws.Range("E37").Value = 'A-nother-Pinsky3^1bn^'
